$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1").Value = 0.98590013770136686
$ws.Range("AC1").Value = 0.78865891421324807
$ws.Range("AR1").Value = 0.80320721207718693
$ws.Range("BD1").Value = 0.93346550727503597
$ws.Range("AS3").Value = 0.7251791499238196
$ws.Range("C4").Value = 0.9083605341823463
$ws.Range("AL4").Value = 0.88202560516350081
$ws.Range("B5").Value = 0.79996188694999137
$ws.Range("BH5").Value = 0.99281978892381229
$ws.Range("H6").Value = 0.94411554051982138
$ws.Range("S7").Value = 0.89910690299524765
$ws.Range("Y7").Value = 0.8188674106293452
$ws.Range("Z8").Value = 0.79471793926032752
$ws.Range("AI8").Value = 0.62499328793269138
$ws.Range("AW8").Value = 0.99125806109847048
$ws.Range("Q9").Value = 0.8996017315776863
$ws.Range("S9").Value = 0.9317695031680171
$ws.Range("L10").Value = 0.71418186272317319
$ws.Range("AO10").Value = 0.99054327089457295
$ws.Range("J11").Value = 0.85849259101038178
$ws.Range("L11").Value = 0.92940802725473959
$ws.Range("AR11").Value = 0.90150707520963214
$ws.Range("AA12").Value = 0.74191197011896293
$ws.Range("BD12").Value = 0.70906727672239378
$ws.Range("Y13").Value = 0.96400903027128404
$ws.Range("AH13").Value = 0.80031952253122474
$ws.Range("I14").Value = 0.92768805610345006
$ws.Range("AC14").Value = 0.72021872511313101
$ws.Range("AY14").Value = 0.9895868498583964
$ws.Range("AR16").Value = 0.6953428764332733
$ws.Range("P18").Value = 0.97391847499783291
$ws.Range("Q18").Value = 0.98773297953041994
$ws.Range("AO19").Value = 0.72811343760248426
$ws.Range("BP19").Value = 0.70293921248844393
$ws.Range("U20").Value = 0.93856198325321216
$ws.Range("C21").Value = 0.84011664797322305
$ws.Range("R21").Value = 0.88074483029841877
$ws.Range("AF21").Value = 0.75185573560689534
$ws.Range("BK21").Value = 0.99773992059205652
$ws.Range("K22").Value = 0.70183477780704939
$ws.Range("R22").Value = 0.92008096902426217
$ws.Range("BN22").Value = 0.86950294321935373
$ws.Range("BO23").Value = 0.95937274393494987
$ws.Range("AB24").Value = 0.90588965290619883
$ws.Range("BG24").Value = 0.70545886995342144
$ws.Range("V25").Value = 0.70559958785189791
$ws.Range("BD26").Value = 0.72789033889242982
$ws.Range("Q27").Value = 0.53683080853973508
$ws.Range("AB27").Value = 0.93073574131093795
$ws.Range("AL27").Value = 0.96919112088875026
$ws.Range("M28").Value = 0.69178780472669832
$ws.Range("BN29").Value = 0.9099990999453178
$ws.Range("J30").Value = 0.96311546545259374
$ws.Range("BI30").Value = 0.82165025811708059
$ws.Range("BA31").Value = 0.90471848335343208
$ws.Range("AG32").Value = 0.78273072518408537
$ws.Range("AB33").Value = 0.81770818937657319
$ws.Range("BF33").Value = 0.90665277419235157
$ws.Range("Q34").Value = 0.90394755509541203
$ws.Range("AR34").Value = 0.82484121683196743
$ws.Range("D36").Value = 0.93707519355146207
$ws.Range("BK36").Value = 0.65452722117628337
$ws.Range("AV37").Value = 0.83050489640962621
$ws.Range("BO37").Value = 0.78472427900911845
$ws.Range("F38").Value = 0.97550060328032973
$ws.Range("AC38").Value = 0.72101874531281118
$ws.Range("AH39").Value = 0.9077990310215891
$ws.Range("AL39").Value = 0.95690874103580592
$ws.Range("AR39").Value = 0.79251655120820552
$ws.Range("P40").Value = 0.76734553021816376
$ws.Range("N41").Value = 0.83166266072077755
$ws.Range("O41").Value = 0.59416983888959729
$ws.Range("AI41").Value = 0.81132741509560713
$ws.Range("Q42").Value = 0.56940644676816332
$ws.Range("U42").Value = 0.90619916101754394
$ws.Range("AN42").Value = 0.80814520187058148
$ws.Range("AX42").Value = 0.93246906708508903
$ws.Range("AF43").Value = 0.76909905494894071
$ws.Range("AR43").Value = 0.67203540174082388
$ws.Range("E44").Value = 0.84110162148413192
$ws.Range("AE44").Value = 0.68015304666942222
$ws.Range("D45").Value = 0.72364861745959186
$ws.Range("BB45").Value = 0.93108023689540387
$ws.Range("AB46").Value = 0.60765532548727719
$ws.Range("BL46").Value = 0.89966614436068693
$ws.Range("BP46").Value = 0.60095980290582651
$ws.Range("W47").Value = 0.86147299450892656
$ws.Range("AT47").Value = 0.81524607691360962
$ws.Range("AV47").Value = 0.97028346243286379
$ws.Range("BF47").Value = 0.55685479638284008
$ws.Range("BG47").Value = 0.93032234663834967
$ws.Range("D48").Value = 0.90965971526289602
$ws.Range("AH48").Value = 0.70202691853135513
$ws.Range("BM48").Value = 0.74025271452084285
$ws.Range("D49").Value = 0.83431418107531896
$ws.Range("AY49").Value = 0.88990321827645724
$ws.Range("AY50").Value = 0.86840660098346123
$ws.Range("AH51").Value = 0.73982469301551723
$ws.Range("BK52").Value = 0.60209306375486704
$ws.Range("I53").Value = 0.62792833840542839
$ws.Range("Z53").Value = 0.89216678358433565
$ws.Range("AB53").Value = 0.93676942111862682
$ws.Range("BF53").Value = 0.96750833346842913
$ws.Range("BG54").Value = 0.92877938856008668
$ws.Range("AC55").Value = 0.68788100557972975
$ws.Range("AR55").Value = 0.7087014489485024
$ws.Range("BI55").Value = 0.95165872094936432
$ws.Range("Y56").Value = 0.82525102572713993
$ws.Range("AM56").Value = 0.6477256551130921
$ws.Range("BG57").Value = 0.96584227388329502
$ws.Range("J58").Value = 0.94228433802979339
$ws.Range("AK59").Value = 0.8111232093822861
$ws.Range("AZ59").Value = 0.9018340253940349
$ws.Range("AS60").Value = 0.85476846291889452
$ws.Range("BL60").Value = 0.91443443345594821
$ws.Range("G61").Value = 0.5957152641458745
$ws.Range("BG61").Value = 0.74696125661708823
$ws.Range("BL61").Value = 0.94762779053566315
$ws.Range("B62").Value = 0.80557186775306744
$ws.Range("O62").Value = 0.69101706406866237
$ws.Range("P62").Value = 0.81595119405654626
$ws.Range("R63").Value = 0.68979336601706875
$ws.Range("AK63").Value = 0.8529284186576076
$ws.Range("BE63").Value = 0.78772977222881768
$ws.Range("BI63").Value = 0.91854360165553506
$ws.Range("BJ63").Value = 0.92199262919999758
$ws.Range("BO63").Value = 0.50889840258806351
$ws.Range("B64").Value = 0.83431931677381299
$ws.Range("BL65").Value = 0.90057468469870672
$ws.Range("BO65").Value = 0.72047411346250301
$ws.Range("B66").Value = 0.85040638356827092
$ws.Range("Y66").Value = 0.71470999761836507
$ws.Range("AV66").Value = 0.91973304504114206
$ws.Range("AG67").Value = 0.8182855878870321
$ws.Range("T68").Value = 0.68193074152656896
$ws.Range("BO68").Value = 0.97074725545670493
